# "saco acentos de los TCs" - update the QA claim numbers (NroSiniestro) in
# rows 2-4 of Hoja1. The leading apostrophe forces Excel to keep/store the
# values as text (preserving leading zeros and any trailing spaces) instead
# of silently converting them to numbers.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = "'0420194406695 "
$ws.Cells.Item(3, 2).Value = "'1220194200662"
$ws.Cells.Item(4, 2).Value = "'1120194100405"

# Update the active selection shown when the sheet is reopened.
$ws.Range("B4").Select()
